$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 <- source row 8 (per rotation described by the diff)
$ws.Range("D4").Value = 44424
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 18000
$ws.Range("N4").Value = '$/caja 15 kilos'
$ws.Range("P4").Value = 1200
$ws.Range("Q4").Value = 15

# Row 5 <- source row 9 (per rotation described by the diff)
$ws.Range("D5").Value = 44424
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = '$/caja 15 kilos'
$ws.Range("P5").Value = 800
$ws.Range("Q5").Value = 15

# Row 6 <- source row 4 (per rotation described by the diff)
$ws.Range("D6").Value = 44536
$ws.Range("J6").Value = 87
$ws.Range("K6").Value = 22000
$ws.Range("L6").Value = 22000
$ws.Range("M6").Value = 22000
$ws.Range("P6").Value = 1222

# Row 7 <- source row 5 (per rotation described by the diff)
$ws.Range("D7").Value = 44536
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 20000
$ws.Range("P7").Value = 1111

# Row 8 <- source row 10 (per rotation described by the diff)
$ws.Range("D8").Value = 44235
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 14000
$ws.Range("N8").Value = '$/bandeja 18 kilos'
$ws.Range("P8").Value = 778
$ws.Range("Q8").Value = 18

# Row 9 <- source row 11 (per rotation described by the diff)
$ws.Range("D9").Value = 44235
$ws.Range("J9").Value = 70
$ws.Range("N9").Value = '$/bandeja 18 kilos'
$ws.Range("P9").Value = 667
$ws.Range("Q9").Value = 18

# Row 10 <- source row 12 (per rotation described by the diff)
$ws.Range("I10").Value = 'Tercera'
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("P10").Value = 556

# Row 11 <- source row 6 (per rotation described by the diff)
$ws.Range("D11").Value = 44242
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = 13000
$ws.Range("P11").Value = 722

# Row 12 <- source row 7 (per rotation described by the diff)
$ws.Range("D12").Value = 44242
$ws.Range("I12").Value = 'Segunda'
$ws.Range("J12").Value = 50
